# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at the top of the data block
# (row 11), pushing all existing records (rows 11..62) down by one row
# (to rows 12..63). This mirrors the upstream append-newest-record-first
# pattern seen in the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 11; Excel shifts rows 11:62 down to 12:63
# and carries the column D date-number style forward onto the new row.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44561
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101004
$ws.Range("J11").Value = "Frambuesa"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("Q11").Value = "`$/bandeja 2 kilos"
$ws.Range("R11").Value = "Provincia de Linares"
$ws.Range("S11").Value = 4000
$ws.Range("T11").Value = 2
